$d = $word.ActiveDocument
$sel = $word.Selection
Write-Host "Selection StoryType:" $sel.StoryType
Write-Host "Selection Text:" $sel.Text
$view = $word.ActiveWindow.ActivePane.View
Write-Host "View.SeekView:" $view.SeekView
$view.SeekView = 9
Write-Host "View.SeekView after:" $view.SeekView
Write-Host "Selection StoryType after seek:" $sel.StoryType
Write-Host "Selection Text after seek:" $sel.Text
